# Weekly data refresh: insert a new price-report row for this market right
# after the header block (new row 24), pushing the existing rows 24:129 down
# to 25:130 (dimension grows from R129 to R130).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(24).Insert()

$ws.Range("A24").Value2 = 2
$ws.Range("B24").Value2 = "Comercializadora del Agro de Limarí"
$ws.Range("C24").Value2 = "Coquimbo"
$ws.Range("D24").Value2 = 45063
$ws.Range("E24").Value2 = 4
$ws.Range("F24").Value2 = 100112030
$ws.Range("G24").Value2 = "Poroto granado"
$ws.Range("H24").Value2 = "Sin especificar"
$ws.Range("I24").Value2 = "Primera"
$ws.Range("J24").Value2 = 800
$ws.Range("K24").Value2 = 23000
$ws.Range("L24").Value2 = 25000
$ws.Range("M24").Value2 = 24000
$ws.Range("N24").Value2 = "$/malla 25 kilos"
$ws.Range("O24").Value2 = "Provincia de Limarí"
$ws.Range("P24").Value2 = 960
$ws.Range("Q24").Value2 = 25
$ws.Range("R24").Value2 = "Hortaliza"
